$d = $word.ActiveDocument

# Locate the paragraph that tells responders "Actions reporter has requested
# you take (do not prompt the reporter):" -- it is being replaced by two new
# prompts that more strongly emphasize not prompting the reporter.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Actions reporter has requested you take \(do not prompt the reporter\):") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $before = $target.Previous()
    $afterBlank1 = $target.Next()
    $afterBlank2 = $afterBlank1.Next()

    # The (empty) paragraph just before the old prompt becomes the new
    # "make the reporter feel safe" question.
    $before.Range.InsertBefore("Is there anything that you can do to make the reporter feel more safe, comfortable, or welcome?")

    # Leave the first blank paragraph after the old prompt untouched, and turn
    # the second blank paragraph into the new "do not prompt" instructions.
    $afterBlank2.Range.InsertBefore("Do NOT prompt the reporter for how to handle the report. If they volunteered actions they want taken, record it here:")

    # Remove the old prompt paragraph (and its paragraph mark) entirely.
    $target.Range.Delete()
}
